$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel stores column width as the value set via ColumnWidth plus a fixed
# padding offset of 5/6 (0.8333...) characters. Subtract that offset so the
# persisted OOXML "width" attribute matches the target values exactly.
$pad = 5 / 6

$ws.Columns.Item(1).ColumnWidth = 66 - $pad
$ws.Columns.Item(2).ColumnWidth = 41 - $pad
$ws.Columns.Item(3).ColumnWidth = 9 - $pad
$ws.Columns.Item(4).ColumnWidth = 49 - $pad
$ws.Columns.Item(5).ColumnWidth = 51 - $pad
$ws.Columns.Item(6).ColumnWidth = 39 - $pad
$ws.Columns.Item(7).ColumnWidth = 41 - $pad
$ws.Columns.Item(8).ColumnWidth = 36 - $pad
$ws.Columns.Item(9).ColumnWidth = 38 - $pad
$ws.Columns.Item(10).ColumnWidth = 44 - $pad
$ws.Columns.Item(11).ColumnWidth = 46 - $pad
$ws.Columns.Item(12).ColumnWidth = 39 - $pad
$ws.Columns.Item(13).ColumnWidth = 41 - $pad

# Header row renames
$ws.Range("B1").Value = "div_testRunDetails_internalRoleCellName"
$ws.Range("D1").Value = "link_testProjectNavigation_internalRoleLinkName"
$ws.Range("E1").Value = "link_testProjectNavigation_internalRoleLinkName_1"
$ws.Range("F1").Value = "link_testProjectNavigation_project_id"
$ws.Range("G1").Value = "link_testProjectNavigation_project_id_1"
$ws.Range("H1").Value = "link_testProjectNavigation_team_id"
$ws.Range("I1").Value = "link_testProjectNavigation_team_id_1"
$ws.Range("J1").Value = "link_testProjectNavigation_test_project_id"
$ws.Range("K1").Value = "link_testProjectNavigation_test_project_id_1"
$ws.Range("L1").Value = "link_testProjectNavigation_trNthChild"
$ws.Range("M1").Value = "link_testProjectNavigation_trNthChild_1"

# Data row update
$ws.Range("A2").Value = "Data Files/AI-Generated/Common/scheduleAndRunTestSuite-test-data"
